$wb = $excel.ActiveWorkbook

# --- 1. Update selection state on "test_index_names" (3rd sheet) -------------
# Old state had this sheet active with a single-cell selection (D2); the new
# state selects columns A:C (full-height) and is no longer the active tab.
$wsIndexNames = $wb.Worksheets.Item(3)
$wsIndexNames.Columns("A:C").Select() | Out-Null

# --- 2. Insert the new "drop_dups" worksheet --------------------------------
# It is placed right after "test_index_names" and before "test_converters",
# i.e. before the (current) 4th sheet.
$wsDropDups = $wb.Worksheets.Add($wb.Worksheets.Item(4))
$wsDropDups.Name = "drop_dups"

# Header row
$wsDropDups.Range("A1").Value = "first_name"
$wsDropDups.Range("B1").Value = "last_name"
$wsDropDups.Range("C1").Value = "wwid"

# Row 2
$wsDropDups.Range("A2").Value = "jonathna"
$wsDropDups.Range("B2").Value = "smith"
$wsDropDups.Range("C2").Value = 1732524

# Row 4 (filled before row 3/5 so shared-string insertion order matches)
$wsDropDups.Range("A4").Value = "austin"
$wsDropDups.Range("B4").Value = "ritter"
$wsDropDups.Range("C4").Value = 423134

# Row 5 - a duplicate-like row that should be dropped
$wsDropDups.Range("A5").Value = "I should"
$wsDropDups.Range("B5").Value = "be removed"
$wsDropDups.Range("C5").Value = 43243

# Row 3 - a look-alike row (trailing space) that should stay
$wsDropDups.Range("A3").Value = "I should "
$wsDropDups.Range("B3").Value = "stay"
$wsDropDups.Range("C3").Value = 43243

# Column widths (best-fit in the source workbook; inputs chosen so the
# engine's internal pixel-rounding lands as close as possible to the
# target character widths of 10.5703125 / 11.7109375 / 8)
$wsDropDups.Columns.Item(1).ColumnWidth = 9.59
$wsDropDups.Columns.Item(2).ColumnWidth = 10.75
$wsDropDups.Columns.Item(3).ColumnWidth = 7.09

# Selection / active cell for the new sheet
$wsDropDups.Range("F8").Select() | Out-Null

# Make the new sheet the active tab (also marks it tabSelected on save)
$wsDropDups.Activate() | Out-Null

Write-Output "drop_dups sheet inserted"
